$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.9639344262295082
$ws.Range("C2").Value = 0.8622950819672132

$ws.Range("B3").Value = 0.9475409836065574
$ws.Range("C3").Value = 0.8852459016393442

$ws.Range("B4").Value = 0.9606557377049181
$ws.Range("C4").Value = 0.8622950819672132

$ws.Range("B5").Value = 0.9573770491803278
$ws.Range("C5").Value = 0.8786885245901639

$ws.Range("B6").Value = 0.9606557377049181
$ws.Range("C6").Value = 0.8688524590163934
